$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ground-truth detection points for rows 3-22 with new data values
$ws.Range("A3").Value = 6186286
$ws.Range("B3").Value = 19
$ws.Range("C3").Value = 41
$ws.Range("D3").Value = 12

$ws.Range("A4").Value = 6242975
$ws.Range("B4").Value = 27
$ws.Range("C4").Value = 59
$ws.Range("D4").Value = 17

$ws.Range("A5").Value = 6254655
$ws.Range("B5").Value = 32
$ws.Range("C5").Value = 43
$ws.Range("D5").Value = 17

$ws.Range("A6").Value = 6293323
$ws.Range("B6").Value = 34
$ws.Range("C6").Value = 44
$ws.Range("D6").Value = 27

$ws.Range("A7").Value = 6358896
$ws.Range("B7").Value = 22
$ws.Range("C7").Value = 44
$ws.Range("D7").Value = 27

$ws.Range("A8").Value = 6393163
$ws.Range("B8").Value = 21
$ws.Range("C8").Value = 49
$ws.Range("D8").Value = 21

$ws.Range("A9").Value = 6717781
$ws.Range("B9").Value = 21
$ws.Range("C9").Value = 52
$ws.Range("D9").Value = 26

$ws.Range("A10").Value = 9189512
$ws.Range("B10").Value = 37
$ws.Range("C10").Value = 58
$ws.Range("D10").Value = 26

$ws.Range("A11").Value = 9185492
$ws.Range("B11").Value = 32
$ws.Range("C11").Value = 36
$ws.Range("D11").Value = 22

$ws.Range("A12").Value = 9157873
$ws.Range("B12").Value = 37
$ws.Range("C12").Value = 48
$ws.Range("D12").Value = 18

$ws.Range("A13").Value = 9156077
$ws.Range("B13").Value = 44
$ws.Range("C13").Value = 56
$ws.Range("D13").Value = 19

$ws.Range("A14").Value = 9155692
$ws.Range("B14").Value = 35
$ws.Range("C14").Value = 54
$ws.Range("D14").Value = 25

$ws.Range("A15").Value = 9118221
$ws.Range("B15").Value = 35
$ws.Range("C15").Value = 52
$ws.Range("D15").Value = 25

$ws.Range("A16").Value = 9109922
$ws.Range("B16").Value = 34
$ws.Range("C16").Value = 34
$ws.Range("D16").Value = 29

$ws.Range("A17").Value = 9096172
$ws.Range("B17").Value = 28
$ws.Range("C17").Value = 57
$ws.Range("D17").Value = 24

$ws.Range("A18").Value = 9073089
$ws.Range("B18").Value = 27
$ws.Range("C18").Value = 41
$ws.Range("D18").Value = 26

$ws.Range("A19").Value = 9073037
$ws.Range("B19").Value = 37
$ws.Range("C19").Value = 44
$ws.Range("D19").Value = 19

$ws.Range("A20").Value = 9065075
$ws.Range("B20").Value = 34
$ws.Range("C20").Value = 45
$ws.Range("D20").Value = 25

$ws.Range("A21").Value = 9064027
$ws.Range("B21").Value = 30
$ws.Range("C21").Value = 44
$ws.Range("D21").Value = 21

$ws.Range("A22").Value = 9058432
$ws.Range("B22").Value = 24
$ws.Range("C22").Value = 51
$ws.Range("D22").Value = 21

# Rows 23-32 no longer have ground truth data - clear their contents
$ws.Range("A23:D28").ClearContents()
$ws.Range("A29:D32").ClearContents()

# Move the active selection to C13 as recorded in the saved view state
$ws.Range("C13").Select()
